$wb = $excel.ActiveWorkbook

# Add the new "InvalidLogin" worksheet after the last existing sheet (ValidLogin)
# so it becomes the 3rd / last tab, matching the target layout.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "InvalidLogin"

# Populate the data-driven test rows for invalid login credentials.
$newSheet.Range("A1").Value = "username"
$newSheet.Range("B1").Value = "password"
$newSheet.Range("A2").Value = "abc"
$newSheet.Range("B2").Value = "xyz"

# Match the saved selection/active cell on the new sheet.
$newSheet.Range("B2").Select() | Out-Null
